$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values for rows 1-21 (row 1 = dates, rows 2-21 = forecast data), columns B-L
$data = @(
    @(45727, 45734, 45741, 45748, 45755, 45762, 45769, 45776, 45783, 45790, 45797),
    @(48.3, 45.1, 42.1, 39.4, 37.2, 35.6, 34.4, 33.8, 33.6, 33.9, 34.6),
    @(65.59999999999999, 63, 59.6, 55.6, 51, 46.3, 41.6, 37.2, 33.3, 30.1, 27.8),
    @(30.2, 26.3, 22.6, 19.3, 16.7, 14.7, 13.3, 12.5, 12.3, 12.8, 14),
    @(47.7, 44.2, 41, 38.2, 36, 34.2, 32.9, 32.1, 31.8, 31.9, 32.4),
    @(66.2, 64.2, 62.1, 60, 58.2, 56.9, 56, 55.7, 55.7, 56.1, 56.6),
    @(42.2, 39.9, 38.3, 37.4, 37.5, 38.2, 39.6, 41.5, 43.7, 46.1, 48.4),
    @(52.2, 48.7, 45.5, 42.7, 40.6, 39.1, 38.1, 37.7, 37.7, 38.2, 39),
    @(39, 36, 33.7, 32.1, 31.2, 31, 31.4, 32.3, 33.7, 35.5, 37.5),
    @(36.9, 34.4, 32.7, 31.9, 32, 32.7, 34.1, 36, 38.2, 40.6, 43.1),
    @(52.3, 48.7, 45.5, 42.9, 41, 39.7, 38.8, 38.6, 38.8, 39.4, 40.4),
    @(56.6, 53, 49.5, 46.3, 43.5, 41.2, 39.2, 37.8, 36.7, 36.1, 36),
    @(59.1, 56.4, 53.7, 50.9, 48.3, 46, 43.8, 42, 40.6, 39.6, 39.1),
    @(43.7, 41, 38.4, 36.2, 34.4, 33.1, 32.3, 32, 32.1, 32.6, 33.4),
    @(37.9, 34, 30.6, 27.8, 25.6, 23.9, 22.7, 21.8, 21.2, 21, 21.2),
    @(48.6, 44.8, 40.8, 37, 33.4, 30.1, 27.2, 24.8, 23, 21.9, 21.4),
    @(55, 50.9, 46.7, 42.6, 38.7, 35.3, 32.3, 29.9, 28.1, 27, 26.7),
    @(42.3, 38.6, 35.3, 32.4, 30.1, 28.3, 27.1, 26.4, 26.3, 26.7, 27.5),
    @(53.9, 51.3, 48.7, 46.4, 44.5, 43.2, 42.6, 42.7, 43.2, 44.2, 45.4),
    @(46.8, 43.3, 40.2, 37.7, 36, 35.1, 35.1, 35.8, 37.2, 39, 41.2),
    @(42, 38.3, 34.7, 31.5, 29, 25.6, 25.7, 25.1, 25, 25.6, 26.6)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 1
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Count; $j++) {
        $colNum = $j + 2
        $ws.Cells.Item($rowNum, $colNum).Value = $rowVals[$j]
    }
}
